# Fruta / hortaliza, semanal
# Update weekly price/volume/origin data for Níspero (Vega Modelo de Temuco)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44496
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 28000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia de Quillota"
$ws.Range("S2").Value = 2800
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44503
$ws.Range("M3").Value = 50
$ws.Range("R3").Value = "Provincia de Quillota"

# Row 4
$ws.Range("D4").Value = 44519
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 28000
$ws.Range("S4").Value = 2800

# Row 6
$ws.Range("D6").Value = 44488
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("S6").Value = 2400

# Row 8
$ws.Range("D8").Value = 44511
$ws.Range("M8").Value = 45
$ws.Range("R8").Value = "Provincia de Los Andes"

# Row 9
$ws.Range("D9").Value = 44511
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 3200
$ws.Range("O9").Value = 3200
$ws.Range("P9").Value = 3200
$ws.Range("S9").Value = 320

# Row 10
$ws.Range("D10").Value = 44515
$ws.Range("M10").Value = 80
$ws.Range("R10").Value = "Provincia de Los Andes"

# Row 11
$ws.Range("D11").Value = 44466
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 11000
$ws.Range("Q11").Value = "$/bandeja 5 kilos"
$ws.Range("R11").Value = "La Ligua"
$ws.Range("S11").Value = 2200
$ws.Range("T11").Value = 5
